$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "EntryTime"
$ws.Range("D1").Value = "ExitTime"

$ws.Range("C2").Value = "12:38:05"
$ws.Range("D2").Value = "12:43:32"
